$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 181 (this shifts the existing rows 181-193
# down to 182-194 and extends the sheet dimension to A1:R194, matching the
# other rows' formatting/styles automatically).
$ws.Rows.Item(181).Insert()

# Populate the newly inserted row 181 with the new record.
$ws.Range("A181").Value = 11
$ws.Range("B181").Value = 'Vega Monumental Concepción'
$ws.Range("C181").Value = 'Bíobío'
$ws.Range("D181").Value = (Get-Date -Year 2023 -Month 10 -Day 13 -Hour 0 -Minute 0 -Second 0)
$ws.Range("E181").Value = 8
$ws.Range("F181").Value = 100112001
$ws.Range("G181").Value = 'Berenjena'
$ws.Range("H181").Value = 'Sin especificar'
$ws.Range("I181").Value = 'Primera'
$ws.Range("J181").Value = 100
$ws.Range("K181").Value = 8000
$ws.Range("L181").Value = 8000
$ws.Range("M181").Value = 8000
$ws.Range("N181").Value = '$/caja 50 unidades'
$ws.Range("O181").Value = 'Región de Arica y Parinacota'
$ws.Range("P181").Value = 160
$ws.Range("Q181").Value = 50
$ws.Range("R181").Value = 'Hortaliza'
